# Update the "Förändrad" date column (C) for rows 2-23 from 2023-09-16 (45185)
# to 2023-10-05 (45204), as described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45185) {
        $cell.Value2 = 45204
    }
}
